$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "sin Low"
$ws.Range("F1").Value = "RAM Megas"
$ws.Range("F1").HorizontalAlignment = -4108

$ws.Range("E2").Value = 300
$ws.Range("E3").Value = 500
$ws.Range("E4").Value = 700

$ws.Range("F2").Value = 1200
$ws.Range("F3").Value = 1740
$ws.Range("F4").Value = 2070

$ws.Range("C4").Select()
